$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.025.10"
$ws.Range("E2").Value = "  +7.97%  "
$ws.Range("D3").Value = "1.822.55"
$ws.Range("E3").Value = "  +5.36%  "
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'246.29"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4930"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("D8").Value = "'44.47"
$ws.Range("E8").Value = "  +7.46%  "
$ws.Range("D9").Value = "'0.2768"
$ws.Range("E9").Value = "  +6.64%  "
$ws.Range("D10").Value = "'0.06383"
$ws.Range("E10").Value = "  +3.37%  "
$ws.Range("D11").Value = "1.820.24"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("D12").Value = "'16.63"
$ws.Range("D13").Value = "'0.07062"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").Value = "'0.6430"
$ws.Range("E14").Value = "  +6.57%  "
$ws.Range("D15").Value = "'84.10"
$ws.Range("E15").Value = "  +9.21%  "
$ws.Range("D16").Value = "'4.693"
$ws.Range("E16").Value = "  +5.20%  "
$ws.Range("D17").Value = "29.030.69"
$ws.Range("E17").Value = "  +8.04%  "
$ws.Range("D18").Value = "'1.0000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'0.000007296"
$ws.Range("E19").Value = "  +2.50%  "
$ws.Range("D20").Value = "'0.9994"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +7.25%  "
$ws.Range("D22").Value = "2.053.85"
$ws.Range("E22").Value = "  +5.23%  "
$ws.Range("D23").Value = "'4.541"
$ws.Range("E23").Value = "  +3.45%  "
$ws.Range("D24").Value = "'8.841"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").Value = "'5.367"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").Value = "'143.47"
$ws.Range("E26").Value = "  +2.62%  "
$ws.Range("D27").Value = "'129.65"
$ws.Range("E27").Value = "  +21.42%  "
$ws.Range("D29").Value = "'1.882"
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D30").Value = "'1.398"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").Value = "'4.124"
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "'0.08359"
$ws.Range("E32").Value = "  +5.63%  "
$ws.Range("D33").Value = "'3.771"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("D34").Value = "'0.04956"
$ws.Range("E34").Value = "  +8.02%  "
$ws.Range("D36").Value = "'2.701"
$ws.Range("E36").Value = "  +4.25%  "
$ws.Range("D37").Value = "'0.6686"
$ws.Range("E37").Value = "  +8.47%  "
$ws.Range("D38").Value = "'2.292"
$ws.Range("E38").Value = "  +15.24%  "
$ws.Range("D39").Value = "'2.687"
$ws.Range("D40").Value = "'0.9478"
$ws.Range("E40").Value = "  +2.55%  "
$ws.Range("D41").Value = "'6.162"
$ws.Range("E41").Value = "  +7.96%  "
$ws.Range("D42").Value = "'0.01581"
$ws.Range("E42").Value = "  +5.63%  "
$ws.Range("D43").Value = "'0.9995"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'101.24"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "'0.4064"
$ws.Range("E45").Value = "  +6.06%  "
$ws.Range("D46").Value = "'7.164"
$ws.Range("E46").Value = "  +5.91%  "
$ws.Range("E47").Value = "  +5.68%  "
$ws.Range("D48").Value = "'0.05541"
$ws.Range("E48").Value = "  +3.30%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.110"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'31.65"
$ws.Range("E50").Value = "  +5.21%  "
$ws.Range("E51").Value = "  +4.52%  "
